# Update countries & provincias Spain
# Applies the 13-Jul-2020 17:14 -> 18:31 data refresh to the "Pais" sheet:
#   - refreshed case totals for a handful of countries (cols B:H)
#   - two pairs + one triple of country-name labels swap position
#     (their row's underlying stats were not re-sorted, so the label text
#     itself moves instead)
#   - the "datos actualizados" timestamp footer string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (row 1, column A) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 18:31"

# --- Country-name label swaps (column A only; row position is unchanged) ---
$ws.Range("A60").Value = "Argelia"
$ws.Range("A61").Value = "Moldavia"

$ws.Range("A91").Value = "Guayana Francesa"
$ws.Range("A92").Value = "Guinea"

$ws.Range("A112").Value = "Libano"
$ws.Range("A113").Value = "Mali"
$ws.Range("A114").Value = "Malaui"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Refreshed statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
function Set-Row($r, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

Set-Row 4   3438025 24030 1518871 1781283 0 89  137871
Set-Row 5   1867841 1665  1213512 582095  0 83  72234
Set-Row 6   904225  24759 569753  310761  0 524 23711
Set-Row 9   317657  2616  286556  24077   0 45  7024
Set-Row 16  243230  169   195106  13157   0 13  34967
Set-Row 19  200323  373   185100  6086    0 3   9137
Set-Row 23  107807  217   71645   27375   0 4   8787
Set-Row 60  19689   494   14019   4652    0 7   1018
Set-Row 61  19439   57    12793   5997    0 7   649
Set-Row 69  13204   30    8353    4498    0 1   353
Set-Row 91  6170    221   3210    2931    0 3   29
Set-Row 92  6141    0     4862    1242    0 0   37
Set-Row 97  4956    31    4086    759     0 0   111
Set-Row 100 3826    23    1374    2259    0 0   193
Set-Row 110 2644    27    1981    652     0 0   11
Set-Row 112 2419    85    1423    960     0 0   36
Set-Row 113 2412    1     1730    561     0 0   121
Set-Row 114 2364    0     557     1769    0 0   38
Set-Row 135 1183    4     1008    165     0 0   10
